$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2517.647
$ws.Range("I32").Value = 1794.5714
$ws.Range("J32").Value = 3023.8
$ws.Range("K32").Value = 1794.5714
$ws.Range("L32").Value = 3023.8
$ws.Range("M32").Value = -1468.5714
$ws.Range("N32").Value = -3675.8
$ws.Range("H51").Value = 14985.286
$ws.Range("J51").Value = 12482.833
$ws.Range("L51").Value = 12482.833
$ws.Range("N51").Value = -13450.833
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 1999
$ws.Range("K62").Value = 1999
$ws.Range("M62").Value = -1375
$ws.Range("H64").Value = 5254.6
$ws.Range("J64").Value = 5567.75
$ws.Range("L64").Value = 5567.75
$ws.Range("N64").Value = -6063.75
$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 1999
$ws.Range("K65").Value = 9995
$ws.Range("M65").Value = -6875
$ws.Range("H67").Value = 5254.6
$ws.Range("J67").Value = 5567.75
$ws.Range("L67").Value = 5567.75
$ws.Range("N67").Value = -7283.75
$ws.Range("H97").Value = 3908.3333
$ws.Range("J97").Value = 3908.3333
$ws.Range("L97").Value = 11724.9999
$ws.Range("N97").Value = -12716.9999
$ws.Range("H129").Value = 2923.625
$ws.Range("I129").Value = 1964.6666
$ws.Range("K129").Value = 5893.9998
$ws.Range("M129").Value = -893.9997999999996
$ws.Range("H138").Value = 1484.7142
$ws.Range("I138").Value = 1599.091
$ws.Range("J138").Value = 1358.9
$ws.Range("K138").Value = 4797.272999999999
$ws.Range("L138").Value = 4076.7
$ws.Range("M138").Value = 342.7270000000008
$ws.Range("N138").Value = -14356.7
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3596.8
$ws.Range("I97").Value = 1996
$ws.Range("K97").Value = 1996
$ws.Range("M97").Value = -1500
$ws.Range("H102").Value = 1665
$ws.Range("I102").Value = 1665
$ws.Range("K102").Value = 1665
$ws.Range("M102").Value = -43
$ws.Range("H106").Value = 32666.334
$ws.Range("J106").Value = 32666.334
$ws.Range("L106").Value = 32666.334
$ws.Range("N106").Value = -35190.334
$ws.Range("H122").Value = 4457.375
$ws.Range("I122").Value = 5109.8335
$ws.Range("K122").Value = 15329.5005
$ws.Range("M122").Value = -12879.5005
$ws.Range("H132").Value = 2681.5881
$ws.Range("I132").Value = 2681.5881
$ws.Range("K132").Value = 8044.7643
$ws.Range("M132").Value = -5514.7643
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12345.333
$ws.Range("I86").Value = 16104.125
$ws.Range("J86").Value = 9338.299999999999
$ws.Range("K86").Value = 16104.125
$ws.Range("L86").Value = 9338.299999999999
$ws.Range("M86").Value = -14981.125
$ws.Range("N86").Value = -11584.3
$ws.Range("H89").Value = 12345.333
$ws.Range("I89").Value = 16104.125
$ws.Range("J89").Value = 9338.299999999999
$ws.Range("K89").Value = 80520.625
$ws.Range("L89").Value = 46691.5
$ws.Range("M89").Value = -74904.625
$ws.Range("N89").Value = -57923.5
$ws.Range("H94").Value = 1557.4
$ws.Range("I94").Value = 2837
$ws.Range("K94").Value = 2837
$ws.Range("M94").Value = -2386
$ws.Range("H105").Value = 2998.4
$ws.Range("I105").Value = 2998
$ws.Range("K105").Value = 2998
$ws.Range("M105").Value = -1251
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1918.6154
$ws.Range("I58").Value = 1232.4445
$ws.Range("J58").Value = 3462.5
$ws.Range("K58").Value = 1232.4445
$ws.Range("L58").Value = 3462.5
$ws.Range("M58").Value = -1029.4445
$ws.Range("N58").Value = -3868.5
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2761.125
$ws.Range("I134").Value = 2458.2
$ws.Range("J134").Value = 3266
$ws.Range("K134").Value = 7374.599999999999
$ws.Range("L134").Value = 9798
$ws.Range("M134").Value = -4839.599999999999
$ws.Range("N134").Value = -14868
$ws.Range("H136").Value = 1918.6154
$ws.Range("I136").Value = 1232.4445
$ws.Range("J136").Value = 3462.5
$ws.Range("K136").Value = 3697.3335
$ws.Range("L136").Value = 10387.5
$ws.Range("M136").Value = -1147.3335
$ws.Range("N136").Value = -15487.5
$ws.Range("H24").Value = 562.5
$ws.Range("I24").Value = 650
$ws.Range("J24").Value = 475
$ws.Range("K24").Value = 1950
$ws.Range("L24").Value = 1425
$ws.Range("M24").Value = -1720
$ws.Range("N24").Value = -1885
$ws.Range("H32").Value = 1298
$ws.Range("I32").Value = 1472
$ws.Range("J32").Value = 950
$ws.Range("K32").Value = 4416
$ws.Range("L32").Value = 2850
$ws.Range("M32").Value = -4133
$ws.Range("N32").Value = -3416
$ws.Range("H34").Value = 3510
$ws.Range("I34").Value = 1044
$ws.Range("J34").Value = 3979.7144
$ws.Range("K34").Value = 3132
$ws.Range("L34").Value = 11939.1432
$ws.Range("M34").Value = -3048
$ws.Range("N34").Value = -12107.1432
$ws.Range("H92").Value = 326.08334
$ws.Range("I92").Value = 381
$ws.Range("J92").Value = 161.33333
$ws.Range("K92").Value = 1143
$ws.Range("L92").Value = 483.99999
$ws.Range("M92").Value = 105
$ws.Range("N92").Value = -2979.99999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 2000
$ws.Range("I130").Value = 2000
$ws.Range("K130").Value = 6000
$ws.Range("M130").Value = -980
$ws.Range("H131").Value = 1950
$ws.Range("J131").Value = 2000
$ws.Range("L131").Value = 6000
$ws.Range("N131").Value = -16080
$ws.Range("H11").Value = 24428572
$ws.Range("I11").Value = 26000000
$ws.Range("J11").Value = 20500002
$ws.Range("K11").Value = 26000000
$ws.Range("L11").Value = 20500002
$ws.Range("M11").Value = -25999861
$ws.Range("N11").Value = -20500280
$ws.Range("H21").Value = 7000007
$ws.Range("J21").Value = 7000007
$ws.Range("L21").Value = 7000007
$ws.Range("N21").Value = -7000353
$ws.Range("H30").Value = 7000007
$ws.Range("J30").Value = 7000007
$ws.Range("L30").Value = 7000007
$ws.Range("N30").Value = -7000217
$ws.Range("H80").Value = 2718.1
$ws.Range("I80").Value = 2785.125
$ws.Range("K80").Value = 2785.125
$ws.Range("M80").Value = -1787.125
$ws.Range("H83").Value = 2718.1
$ws.Range("I83").Value = 2785.125
$ws.Range("K83").Value = 13925.625
$ws.Range("M83").Value = -8933.625
$ws.Range("H102").Value = 1128.5
$ws.Range("I102").Value = 985.625
$ws.Range("K102").Value = 985.625
$ws.Range("M102").Value = 636.375
$ws.Range("H126").Value = 20000
$ws.Range("J126").Value = 20000
$ws.Range("L126").Value = 60000
$ws.Range("N126").Value = -64940
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 2500
$ws.Range("I26").Value = 2500
$ws.Range("K26").Value = 2500
$ws.Range("M26").Value = -2205
$ws.Range("H46").Value = 1614.2142
$ws.Range("I46").Value = 1128.6666
$ws.Range("K46").Value = 1128.6666
$ws.Range("M46").Value = -940.6666
$ws.Range("H55").Value = 903.82355
$ws.Range("I55").Value = 283.375
$ws.Range("J55").Value = 1455.3334
$ws.Range("K55").Value = 283.375
$ws.Range("L55").Value = 1455.3334
$ws.Range("M55").Value = -110.375
$ws.Range("N55").Value = -1801.3334
$ws.Range("H68").Value = 4624.5
$ws.Range("J68").Value = 4499.6665
$ws.Range("L68").Value = 4499.6665
$ws.Range("N68").Value = -5997.6665
$ws.Range("H71").Value = 4624.5
$ws.Range("J71").Value = 4499.6665
$ws.Range("L71").Value = 22498.3325
$ws.Range("N71").Value = -29986.3325
$ws.Range("H100").Value = 1669.0667
$ws.Range("I100").Value = 1625.9231
$ws.Range("K100").Value = 1625.9231
$ws.Range("M100").Value = -1084.9231
$ws.Range("H122").Value = 3620.75
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 31577
$ws.Range("I136").Value = 1673.5
$ws.Range("J136").Value = 55499.8
$ws.Range("K136").Value = 5020.5
$ws.Range("L136").Value = 166499.4
$ws.Range("M136").Value = -2470.5
$ws.Range("N136").Value = -171599.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 119994.5
$ws.Range("I29").Value = 119994.5
$ws.Range("K29").Value = 119994.5
$ws.Range("M29").Value = -119704.5
$ws.Range("H69").Value = 72758.5
$ws.Range("J69").Value = 100271
$ws.Range("L69").Value = 100271
$ws.Range("N69").Value = -101769
$ws.Range("H72").Value = 72758.5
$ws.Range("J72").Value = 100271
$ws.Range("L72").Value = 300813
$ws.Range("N72").Value = -308301
$ws.Range("H126").Value = 5101.9165
$ws.Range("I126").Value = 5828.8335
$ws.Range("K126").Value = 17486.5005
$ws.Range("M126").Value = -15016.5005
$ws.Range("H132").Value = 3373.0833
$ws.Range("I132").Value = 2831.889
$ws.Range("J132").Value = 4996.6665
$ws.Range("K132").Value = 8495.667000000001
$ws.Range("L132").Value = 14989.9995
$ws.Range("M132").Value = -5965.667000000001
$ws.Range("N132").Value = -20049.9995
$ws.Range("H136").Value = 13736.723
$ws.Range("I136").Value = 8964.134
$ws.Range("K136").Value = 26892.402
$ws.Range("M136").Value = -24342.402
